# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Refreshes several MSME indicator figures for Latvia with more precise
# (two decimal place) values. Cells hold these figures as text, so a
# leading apostrophe is used to force text storage (matching the
# workbook's existing convention of storing numeric-looking figures as
# strings) instead of letting Excel auto-convert them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro / MSMEs
$ws.Range("B13").Value = "'34.21"
$ws.Range("D13").Value = "'38.41"

# Employment (% of total): Micro / SMEs / MSMEs
$ws.Range("B14").Value = "'25.35"
$ws.Range("C14").Value = "'51.11"
$ws.Range("D14").Value = "'76.46"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B16").Value = "'88.84"
$ws.Range("C16").Value = "'10.92"
$ws.Range("D16").Value = "'99.76"

# Value added to the economy (% of total): Micro / SMEs / MSMEs
$ws.Range("B20").Value = "'15.91"
$ws.Range("C20").Value = "'52.04"
$ws.Range("D20").Value = "'67.95"
